$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New weekly column: AC = 18_05_2021, mirrors AB except rows 8 & 9 which bump by 2
$ws.Range("AC1").Value = "18_05_2021"

$ws.Range("AC2").Value = 1
$ws.Range("AC3").Value = 0
$ws.Range("AC4").Value = 0
$ws.Range("AC5").Value = 7
$ws.Range("AC6").Value = 8
$ws.Range("AC7").Value = 62
$ws.Range("AC8").Value = 210
$ws.Range("AC9").Value = 663
$ws.Range("AC10").Value = 991
$ws.Range("AC11").Value = 561

$ws.Range("AC12").Formula = "=SUM(AC2:AC11)"

$ws.Range("AB12:AC12").Select()
